$wb = $excel.ActiveWorkbook

$buildTimes = @(34128, 34128, 45504, 56880, 68256, 91008, 113760, 136512, 159264, 182016, 204768, 204768, 204768, 204768, 204768)

foreach ($ws in $wb.Worksheets) {
    $ws.Range("G1").Value = "INT_buildTime"
    $ws.Range("G1").Style = $ws.Range("D1").Style
    for ($i = 0; $i -lt $buildTimes.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 7).Value = $buildTimes[$i]
    }
}
